$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.877.59'
$ws.Range('E2').Value = '  -2.31%  '
$ws.Range('D3').Value = '2.451.72'
$ws.Range('E3').Value = '  -3.65%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.24'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.74'
$ws.Range('E6').Value = '  -3.73%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.564'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('E9').Value = '  -1.39%  '
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '4.95'
$ws.Range('E11').Value = '  -5.18%  '
$ws.Range('E12').Value = '  -4.26%  '
$ws.Range('D13').Value = '2.886.02'
$ws.Range('E13').Value = '  -3.63%  '
$ws.Range('D14').Value = '57.805.68'
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.56'
$ws.Range('E15').Value = '  -4.00%  '
$ws.Range('E16').Value = '  -3.31%  '
$ws.Range('D17').Value = '2.455.75'
$ws.Range('E17').Value = '  -3.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.38'
$ws.Range('E18').Value = '  -3.42%  '
$ws.Range('E19').Value = '  -2.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '312.74'
$ws.Range('E20').Value = '  -3.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.11'
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('E23').Value = '  -0.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.403'
$ws.Range('E24').Value = '  -2.13%  '
$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.568.16'
$ws.Range('E26').Value = '  -3.46%  '
$ws.Range('E27').Value = '  -3.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.25'
$ws.Range('E28').Value = '  -3.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '174.49'
$ws.Range('E29').Value = '  +3.27%  '
$ws.Range('E30').Value = '  -3.49%  '
$ws.Range('E31').Value = '  -2.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.14'
$ws.Range('E32').Value = '  -3.70%  '
$ws.Range('E33').Value = '  -7.55%  '
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.81'
$ws.Range('E36').Value = '  -2.77%  '
$ws.Range('E37').Value = '  -7.54%  '
$ws.Range('E38').Value = '  -5.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.28'
$ws.Range('E39').Value = '  -1.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.808'
$ws.Range('E40').Value = '  +2.50%  '
$ws.Range('E41').Value = '  -4.92%  '
$ws.Range('E42').Value = '  -3.06%  '
$ws.Range('E43').Value = '  -3.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.79'
$ws.Range('E44').Value = '  -6.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '257.66'
$ws.Range('E45').Value = '  -8.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '123.32'
$ws.Range('E46').Value = '  -7.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0919'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('E48').Value = '  -3.09%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0211'
$ws.Range('E49').Value = '  -3.22%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.05'
$ws.Range('E50').Value = '  -4.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.25'
$ws.Range('E51').Value = '  -5.60%  '
